# Add skill "Tian Shan Liu Yang Zhang" as a new block of 6 rows
# right before the existing "Increase attack by 10%" level-table block
# (old row 40), i.e. the new block becomes rows 40-45 and everything
# from the old row 40 onward shifts down by 6 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert 6 blank rows at 40..45 (old row 40 -> new row 46).
$ws.Rows("40:45").Insert()

# 2) Header row for the new skill (row 40) - mirrors the layout used by
#    every other skill block (e.g. row 34 "Luo Han Quan").
$ws.Cells.Item(40, 2).Value = "TIAN_SHAN_LIU_YANG_ZHANG"
$ws.Cells.Item(40, 3).Value = "Tian Shan Liu Yang Zhang"
$ws.Cells.Item(40, 4).Formula = "=K40"
$ws.Cells.Item(40, 6).Value = "BEHAVIOR_TYPE"
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 1
$ws.Cells.Item(40, 9).Value = "NONE"
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = "Strike through the end and return back"

# 3) Cost / attribute rows 41-45.
$ws.Cells.Item(41, 6).Value = "ATTACK"
$ws.Cells.Item(41, 7).Value = 20

$ws.Cells.Item(42, 6).Value = "CRITICAL_RATE"
$ws.Cells.Item(42, 7).Value = 0.15

$ws.Cells.Item(43, 6).Value = "CD"
$ws.Cells.Item(43, 7).Value = 4

$ws.Cells.Item(44, 6).Value = "SPEED"
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 9).ClearContents()

$ws.Cells.Item(45, 6).Value = "PIERCE"
$ws.Cells.Item(45, 7).Value = 4
$ws.Cells.Item(45, 9).ClearContents()

# Row 43 has no I-column cell at all in the target layout (unlike 44/45) -
# make sure nothing lingers there.
$ws.Cells.Item(43, 9).ClearContents()
$ws.Cells.Item(41, 9).ClearContents()
$ws.Cells.Item(42, 9).ClearContents()

# 4) Visual formatting: the header row of every skill block carries a thin
#    top border across A:K, and the derived-description cell (column D) is
#    rendered in the same blue tone used elsewhere in the sheet.
$ws.Range("A40:K40").Borders.Item(8).LineStyle = 1
$ws.Cells.Item(40, 4).Font.Color = 12611584

# 5) Keep the sheet's selection / active cell consistent with the edit.
$ws.Range("G45").Select()
